$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.10522383444921957
$ws.Range("A2").Value = -0.053648265261436734
$ws.Range("A3").Value = -0.0089999994391298799
$ws.Range("A4").Value = 0.28399221189407342
$ws.Range("A5").Value = -0.0059999994520252287
$ws.Range("A6").Value = -0.029995200392043841
$ws.Range("A7").Value = -0.019999999343017549
$ws.Range("A8").Value = -0.019999999339887609
$ws.Range("A9").Value = -0.0059999994241382026
$ws.Range("A10").Value = -0.005999999422392932
$ws.Range("A11").Value = -0.0044999994316121672
$ws.Range("A12").Value = 0.045862733032532521
$ws.Range("A13").Value = -0.0095167616315743331
$ws.Range("A14").Value = -0.01199999937222973
$ws.Range("A15").Value = -0.0059999994076953556
$ws.Range("A16").Value = -0.0059999994060708772
$ws.Range("A17").Value = -0.0059999994037731597
$ws.Range("A18").Value = -0.0089999993847671433
$ws.Range("A19").Value = -0.0089999994401117611
$ws.Range("A20").Value = -0.0089999994355132173
$ws.Range("A21").Value = -0.0089999994349270196
$ws.Range("A22").Value = -0.0089999994345424383
$ws.Range("A23").Value = -0.0089999994356402269
$ws.Range("A24").Value = -0.041999999222458229
$ws.Range("A25").Value = -0.041999999218595541
$ws.Range("A26").Value = -0.005999999429242564
$ws.Range("A27").Value = -0.0059999994269763768
$ws.Range("A28").Value = -0.0059999994151525016
$ws.Range("A29").Value = -0.011999999370329917
$ws.Range("A30").Value = -0.019999999316753669
$ws.Range("A31").Value = -0.014999999342283132
$ws.Range("A32").Value = -0.020999999304220474
$ws.Range("A33").Value = -0.0059999993970762944
